# Atualização da Lista de Riscos
# Applies the "Sistema de Rastreamento" risk-list update:
#   - fixes the identification dates for risks 1-4
#   - fills in risks 5 and 6 (rows 7 and 8) with full data
#   - widens column I slightly
#   - updates the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the "Data de Identificação" values for the first four risks ---
$ws.Range("B3").Value = 41842
$ws.Range("B4").Value = 41869
$ws.Range("B5").Value = 41871
$ws.Range("B6").Value = 41886

# --- Risk 5 (row 7): power-supply risk ---
$ws.Range("B7").Value = 41888
# Description (D) must be entered before Name (C) so new shared strings
# are appended to xl/sharedStrings.xml in the expected order.
$ws.Range("D7").Value = "É possível que ocorram problemas diversos com o fornecimento de energia elétrica. Este problema pode afetar não apenas o desenvolvimento, como também a queima de equipamentos diversos."
$ws.Range("C7").Value = "Fornecimento de energia elétrica ininterrupta para todos os equipamentos que serão utilizados pela a equipe de desenvolvimento de software. "
$ws.Range("E7").Value = "D"
$ws.Range("F7").Value = 5
$ws.Range("G7").Value = 0.2
$ws.Range("I7").Value = "Gerente de Configuração e Mudança."
$ws.Range("J7").Value = "Contratar empresa de locação de Gerador Elétrico e que será responsável também pela prestação de serviço de instalação e manutenção do mesmo."
$ws.Rows.Item(7).RowHeight = 75.75

# --- Risk 6 (row 8): equipment-burn-out risk ---
$ws.Range("B8").Value = 41894
$ws.Range("C8").Value = "Queima de equipamentos diversos que serão utilizados pela a equipe de desenvolvimento de software"
$ws.Range("D8").Value = "Pode ser necessário realizarmos o conserto ou substituição de equipamentos que vierem a apresentar defeitos."
$ws.Range("E8").Value = "I"
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.05
$ws.Range("I8").Value = "Gerente de Configuração e Mudança."
$ws.Range("J8").Value = "Contratar empresa de manutenção para os equipamentos elétrico/eletrônicos que serão utilizados por toda a equipe de desenvolvimento de software."
$ws.Rows.Item(8).RowHeight = 63.75

# --- Column I widened to fit the new "Responsável" text ---
$ws.Columns.Item(9).ColumnWidth = 11.29

# --- Update the active selection on the sheet ---
$ws.Range("K8").Select()
